# Auto-generated Excel COM-interop script
# Updates the cryptos list (prices / 1h volume %) and fixes a few
# B/C/D/E row mis-orderings, per the GitHub Actions data refresh commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'45.276.60"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +4.87%  "

# Row 3
$ws.Range("D3").Value = "'2.372.13"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.96%  "

# Row 4
$ws.Range("E4").Value = "  -0.07%  "

# Row 5
$ws.Range("D5").Value = "'109.09"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.86%  "

# Row 6
$ws.Range("D6").Value = "'310.49"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.13%  "

# Row 7
$ws.Range("D7").Value = "'0.631"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.18%  "

# Row 8
$ws.Range("E8").Value = "  -0.21%  "

# Row 9
$ws.Range("E9").Value = "  +0.55%  "

# Row 10
$ws.Range("D10").Value = "'41.25"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.01%  "

# Row 11
$ws.Range("D11").Value = "'0.0919"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.21%  "

# Row 12
$ws.Range("D12").Value = "'8.51"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.10%  "

# Row 13
$ws.Range("E13").Value = "  +1.50%  "

# Row 15
$ws.Range("D15").Value = "'2.731.72"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.84%  "

# Row 16
$ws.Range("D16").Value = "'15.29"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.25%  "

# Row 17
$ws.Range("D17").Value = "'2.368.70"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.84%  "

# Row 18
$ws.Range("D18").Value = "'45.255.31"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +4.35%  "

# Row 19
$ws.Range("D19").Value = "'14.92"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +13.43%  "

# Row 20
$ws.Range("E20").Value = "  -3.17%  "

# Row 21
$ws.Range("E21").Value = "  -0.05%  "

# Row 22
$ws.Range("D22").Value = "'73.35"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.03%  "

# Row 23
$ws.Range("D23").Value = "'3.49"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.16%  "

# Row 24
$ws.Range("D24").Value = "'260.34"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.74%  "

# Row 25
$ws.Range("E25").Value = "  +2.23%  "

# Row 26
$ws.Range("E26").Value = "  -0.31%  "

# Row 27
$ws.Range("D27").Value = "'11.16"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.05%  "

# Row 28
$ws.Range("D28").Value = "'7.31"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -3.65%  "

# Row 29
$ws.Range("D29").Value = "'2.35"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.46%  "

# Row 30
$ws.Range("E30").Value = "  +9.04%  "

# Row 31
$ws.Range("D31").Value = "'22.41"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.85%  "

# Row 32
$ws.Range("D32").Value = "'37.80"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.42%  "

# Row 33
$ws.Range("D33").Value = "'169.28"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.40%  "

# Row 34
$ws.Range("E34").Value = "  +5.15%  "

# Row 35
$ws.Range("E35").Value = "  -0.43%  "

# Row 36
$ws.Range("E36").Value = "  +3.94%  "

# Row 37
$ws.Range("D37").Value = "'4.75"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.43%  "

# Row 38
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").Value = "'3.95"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +4.24%  "

# Row 39
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").Value = "'2.96"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +3.72%  "

# Row 40
$ws.Range("E40").Value = "  -2.38%  "

# Row 41
$ws.Range("D41").Value = "'1.76"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.64%  "

# Row 42
$ws.Range("D42").Value = "'100.04"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -4.44%  "

# Row 43
$ws.Range("D43").Value = "'0.232"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.98%  "

# Row 44
$ws.Range("D44").Value = "'69.60"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -3.05%  "

# Row 45
$ws.Range("D45").Value = "'13.02"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.52%  "

# Row 46
$ws.Range("E46").Value = "  -0.24%  "

# Row 47
$ws.Range("B47").Value = "ordi"
$ws.Range("C47").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D47").Value = "'81.62"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +9.33%  "

# Row 48
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "'1.752.24"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +5.62%  "

# Row 49
$ws.Range("B49").Value = "THORChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D49").Value = "'5.57"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +4.46%  "

# Row 50
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "'112.17"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.66%  "

# Row 51
$ws.Range("E51").Value = "  +2.84%  "
